$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.882.76"
$ws.Range("E2").Value = "  -2.15%  "
$ws.Range("D3").Value = "2.299.28"
$ws.Range("E3").Value = "  -5.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.42%  "
$ws.Range("E8").Value = "  -2.85%  "
$ws.Range("D9").Value = "2.297.49"
$ws.Range("E9").Value = "  -5.15%  "
$ws.Range("E10").Value = "  -3.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.43%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("E13").Value = "  -4.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.39%  "
$ws.Range("D15").Value = "2.706.51"
$ws.Range("E15").Value = "  -5.29%  "
$ws.Range("D16").Value = "58.847.36"
$ws.Range("E16").Value = "  -2.09%  "
$ws.Range("E17").Value = "  -3.25%  "
$ws.Range("D18").Value = "2.305.48"
$ws.Range("E18").Value = "  -4.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "315.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.48%  "
$ws.Range("E22").Value = "  -3.95%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.24%  "
$ws.Range("E25").Value = "  -2.76%  "
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  -6.34%  "
$ws.Range("E28").Value = "  -7.80%  "
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("E31").Value = "  -5.44%  "
$ws.Range("E32").Value = "  +4.13%  "
$ws.Range("E33").Value = "  -4.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.92%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -6.12%  "
$ws.Range("E39").Value = "  -5.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("E41").Value = "  -5.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "303.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "140.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("E44").Value = "  -5.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0954"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0502"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.559"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.18%  "
$ws.Range("E49").Value = "  -3.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.42%  "
$ws.Range("E51").Value = "  -0.29%  "
